$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The first table (rows 2-10) ends with the "Lämna en recension" row at
# row 10 (thick bottom border). A new test case ("Radera användare") is
# being inserted as the new row 10, pushing the old row 10 content down
# into row 11 (which was previously a blank spacer row before the
# "Admin" table header at row 12).

# 1) Insert a blank row above the current row 10. This shifts the old
#    row 10 (and everything below it, including the blank spacer row)
#    down by one.
$ws.Rows("10:10").Insert(-4121)

# 2) The row that used to be the blank spacer (pushed down to row 12)
#    is no longer needed, since the old row 10 now occupies row 11 and
#    naturally borders the "Admin" header row. Remove it so row numbers
#    below stay exactly as they were (12-17).
$ws.Rows("12:12").Delete()

# 3) The newly inserted row 10 has no formatting yet - copy the format
#    from row 9 (the row directly above, same "normal" non-thick-bottom
#    style used by every non-final row of the table).
$ws.Range("A9:G9").Copy()
$ws.Range("A10:G10").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# 4) Fill in the new test case data in row 10.
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Radera användare"
$ws.Range("C10").Value = 'Välj "Användarsida" och tryck på "Ta bort konto"'
$ws.Range("D10").Value = "Kontot tas bort och användaren loggas ut"
$ws.Range("E10").Value = "Kontot är raderat"

# 5) Update the active selection to match the authored change.
$ws.Range("E10").Select()
